# CRS.xlsx - "Doc: updated according to coach comments"
# The author trimmed two rows from Sheet1 (a stray duplicate row inside the
# Client requirements block, and a leftover blank row near the Supplier
# block) and reworded several of the Client requirement rows (cart /
# checkout / client-id wording).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop one row out of the Client block (18) and the stray blank row that
# used to sit right before "CRS_supplier_006" (31) - this takes the sheet
# from 53 rows down to 51, matching the new dimension.
$ws.Rows(18).Delete()
$ws.Rows(31).Delete()

# Re-populate the Client requirement rows with the revised wording / ids.
$ws.Range("A18").Value = 'CRS_Client_002'
$ws.Range("B18").Value = 'client id should consists of the last 6 numbers of national id'
$ws.Range("C18").Value = 'CRS-SIQ_Client_002'

$ws.Range("A19").Value = 'CRS_Client_003'
$ws.Range("B19").Value = 'The client should be able to select products to add one or more to a shopping cart for future purchase consideration.'
$ws.Range("C19").Value = 'CRS-SIQ_Client_003'

$ws.Range("A20").Value = 'CRS_Client_004'
$ws.Range("B20").Value = "The client should be able to view and review the contents of their cart, including product names, quantities, prices, `nand the ability to increment or decrement product count or remove items."
$ws.Range("C20").Value = 'CRS-SIQ_Client_003'

$ws.Range("A21").Value = 'CRS_Client_005'
$ws.Range("B21").Value = 'If the cart is empty, a message "Your cart is empty." should be displayed.'
$ws.Range("C21").Value = 'CRS-SIQ_Client_004'

$ws.Range("A22").Value = 'CRS_Client_006'
$ws.Range("B22").Value = "The client should be able to navigate to the checkout page to review the final list of items and choose the payment method (Cash Only)`nbefore confirming the purchase."
$ws.Range("C22").Value = 'CRS-SIQ_Client_005'

$ws.Range("A23").Value = 'CRS_Client_006'
$ws.Range("B23").Value = 'The client should be able to place an order by confirming the cart contents and completing the checkout process.'
$ws.Range("C23").Value = 'CRS-SIQ_Client_006'

$ws.Range("A24").Value = 'CRS_Client_007'
$ws.Range("B24").Value = 'After placing order confirmation message should appears "your orders placed successfully" and redirect to history of buys'
$ws.Range("C24").Value = 'CRS-SIQ_Client_006'

$ws.Range("A25").Value = 'CRS_Client_008'
$ws.Range("B25").Value = "The client can view a complete history of their previous purchases including (product details, order dates,items purchased, price, `nshipping address and order status)."
$ws.Range("C25").Value = 'CRS-SIQ_Client_007'
